$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

# Update LTP (B) and PREV (C) values for rows 2-25
$values = @{
    2  = @(2822.15, 2887.15)
    3  = @(458.55, 455.15)
    4  = @(1710.15, 1704.05)
    5  = @(7307.6, 7391.2)
    6  = @(237.05, 238.15)
    7  = @(212, 211.15)
    8  = @(47499.8, 47047.15)
    9  = @(649.9, 648.6)
    10 = @(4055, 4091)
    11 = @(154.5, 154.65)
    12 = @(1364.1, 1328)
    13 = @(520.35, 516.05)
    14 = @(1508.95, 1504.75)
    15 = @(686.25, 684)
    16 = @(532.3, 518.7)
    17 = @(1668.55, 1692.65)
    18 = @(273.25, 274.15)
    19 = @(21075, 21030.8)
    20 = @(285.05, 283.8)
    21 = @(614.15, 611.7)
    22 = @(671.6, 679.5)
    23 = @(714.55, 721.95)
    24 = @(323.55, 325.8)
    25 = @(129.2, 130)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}

# Update the active cell selection
$ws.Range("L9").Select()
